# Applies the "Added additional column in read data from excel file" change:
#  - Test Cases sheet: clear the (now unused) per-test-case Result column values
#  - Login_001 / Login_002 sheets: insert a new "Data" column before the
#    Result (now last) column, populate it (plain text, or a hyperlink for
#    the "Navigate to URL" row), rename the enterUserName/enterPassword
#    actions to enterText, and blank out the old Result column values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Test Cases": clear D2:D3 (Result) but keep the bordered cell style
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Range("D2").Value = ""
$wsCases.Range("D3").Value = ""
$wsCases.Range("C2").Copy()
$wsCases.Range("D2").PasteSpecial(-4122)
$wsCases.Range("C3").Copy()
$wsCases.Range("D3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Insert the new "Data" column (before the old "Result" column) on both
# Login_001 and Login_002, and fill in its header + the per-row sample data.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Login_001")
$ws2.Columns("F:F").Insert()
$ws2.Range("F1").Value = "Data"
$ws2.Range("F4").Value = "german"
$ws2.Range("F5").Value = "password"

$ws3 = $wb.Worksheets.Item("Login_002")
$ws3.Columns("F:F").Insert()
$ws3.Range("F1").Value = "Data"
$ws3.Range("F4").Value = "Arabic"
$ws3.Range("F5").Value = "password"

# "Navigate to URL" row (3) gets an actual hyperlink instead of plain text.
$ws2.Hyperlinks.Add($ws2.Range("F3"), "http://demo.silverstripe.org/Security/login") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "http://demo.silverstripe.org/Security/login") | Out-Null

# The username/password steps now drive through a generic "enterText" action.
$ws2.Range("E4").Value = "enterText"
$ws2.Range("E5").Value = "enterText"
$ws2.Range("E6").Value = "click"
$ws3.Range("E4").Value = "enterText"
$ws3.Range("E5").Value = "enterText"
$ws3.Range("E6").Value = "click"

# Blank out the old "Result" column (now G) but keep a bordered style on
# every data row so the cells still round-trip as styled-empty cells.
$r = 2
while ($r -le 8) {
    $ws2.Range("G$r").Value = ""
    $ws2.Range("E$r").Copy()
    $ws2.Range("G$r").PasteSpecial(-4122)
    $ws3.Range("G$r").Value = ""
    $ws3.Range("E$r").Copy()
    $ws3.Range("G$r").PasteSpecial(-4122)
    $r = $r + 1
}

$ws2.Range("A2").Select()
$ws3.Range("E6").Select()
$ws3.Activate()
